# Adds Sheet2 with scratch-notes tables for three demographic models
# (two_epoch, exponential, bottleneck, three_epoch) and repositions the
# active-sheet/selection state to match the authoring session.

$wb = $excel.ActiveWorkbook

# --- clear Sheet1's previous "last touched" selection -----------------
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("I2").Select()

# --- add the new worksheet, placed after Sheet1 -----------------------
$ws = $wb.Worksheets.Add($null, $ws1)
$ws.Name = "Sheet2"

  $ws.Range("A1").Value = "two_epoch"

  $ws.Range("B2").Value = "b_vulgatus"
  $ws.Range("C2").Value = "b_ovatus"
  $ws.Range("D2").Value = "a_putredinis"
  $ws.Range("E2").Value = "b_uniformis"
  $ws.Range("F2").Value = "e_rectale"
  $ws.Range("G2").Value = "p_copri"
  $ws.Range("H2").Value = "o_splanchnicus"
  $ws.Range("I2").Value = "a_shahii"
  $ws.Range("J2").Value = "coprococcus_sp"

  $ws.Range("A3").Value = "Theta"
  $ws.Range("B3").Value = 2026.6041198099999
  $ws.Range("G3").Value = 6851.08948573

  $ws.Range("A4").Value = "Nu_a"
  $ws.Range("B4").Value = 1.49516307
  $ws.Range("G4").Value = 3.9692834100000001

  $ws.Range("A5").Value = "Nu_b"

  $ws.Range("A6").Value = "T_12"
  $ws.Range("B6").Value = 0.44708505999999998
  $ws.Range("G6").Value = 0.39557179999999997

  $ws.Range("A7").Value = "T_23"

  $ws.Range("A8").Value = "N_e"
  $ws.Range("B8").Formula = "=B3/4*1000000000"
  $ws.Range("C8:K8").Formula = "=C3/4*1000000000"

  $ws.Range("A9").Value = "N_a"
  $ws.Range("B9").Formula = "=B8/B4"
  $ws.Range("C9:K9").Formula = "=C8/C4"

  $ws.Range("A10").Value = "T_12 (scaled)"
  $ws.Range("B10").Formula = "=2*B6*B9"
  $ws.Range("C10:K10").Formula = "=2*C6*C9"

  $ws.Range("A11").Value = "T_23 (scaled)"
  $ws.Range("B11").Formula = "=2*B7*B9"
  $ws.Range("C11:K11").Formula = "=2*C7*C9"

  $ws.Range("A13").Value = "exponential"

  $ws.Range("B14").Value = "b_vulgatus"
  $ws.Range("C14").Value = "b_ovatus"
  $ws.Range("D14").Value = "a_putredinis"
  $ws.Range("E14").Value = "b_uniformis"
  $ws.Range("F14").Value = "e_rectale"
  $ws.Range("G14").Value = "p_copri"
  $ws.Range("H14").Value = "o_splanchnicus"
  $ws.Range("I14").Value = "a_shahii"
  $ws.Range("J14").Value = "coprococcus_sp"

  $ws.Range("A15").Value = "Theta"
  $ws.Range("B15").Value = 1939.1988649299999
  $ws.Range("G15").Value = 6651.9430162500003

  $ws.Range("A16").Value = "Nu_a"
  $ws.Range("B16").Value = 1.62764391
  $ws.Range("G16").Value = 5.5207001800000004

  $ws.Range("A17").Value = "Nu_b"
  $ws.Range("B17").Value = "N/A"
  $ws.Range("G17").Value = "N/A"

  $ws.Range("A18").Value = "T_12"
  $ws.Range("B18").Value = 1.09100781
  $ws.Range("G18").Value = 0.65275521000000003

  $ws.Range("A19").Value = "T_23"

  $ws.Range("A20").Value = "N_e"
  $ws.Range("B20").Formula = "=B15/4*1000000000"
  $ws.Range("C20:K20").Formula = "=C15/4*1000000000"

  $ws.Range("A21").Value = "N_a"
  $ws.Range("B21").Formula = "=B20/B16"
  $ws.Range("C21:K21").Formula = "=C20/C16"

  $ws.Range("A22").Value = "T_12 (scaled)"
  $ws.Range("B22").Formula = "=2*B18*B21"
  $ws.Range("C22:K22").Formula = "=2*C18*C21"

  $ws.Range("A23").Value = "T_23 (scaled)"
  $ws.Range("B23").Formula = "=2*B19*B21"
  $ws.Range("C23:K23").Formula = "=2*C19*C21"

  $ws.Range("A25").Value = "bottleneck"

  $ws.Range("B26").Value = "b_vulgatus"
  $ws.Range("C26").Value = "b_ovatus"
  $ws.Range("D26").Value = "a_putredinis"
  $ws.Range("E26").Value = "b_uniformis"
  $ws.Range("F26").Value = "e_rectale"
  $ws.Range("G26").Value = "p_copri"
  $ws.Range("H26").Value = "o_splanchnicus"
  $ws.Range("I26").Value = "a_shahii"
  $ws.Range("J26").Value = "coprococcus_sp"

  $ws.Range("A27").Value = "Theta"
  $ws.Range("B27").Value = 1931.3930001599999
  $ws.Range("G27").Value = 6649.6250826400001

  $ws.Range("A28").Value = "Nu_a"
  $ws.Range("B28").Value = 1.07053664
  $ws.Range("G28").Value = 1.00120516

  $ws.Range("A29").Value = "Nu_b"
  $ws.Range("B29").Value = 1.6312748500000001
  $ws.Range("G29").Value = 5.5203546100000001

  $ws.Range("A30").Value = "T_12"
  $ws.Range("B30").Value = 0.98098777000000004
  $ws.Range("G30").Value = 0.65317188000000004

  $ws.Range("A31").Value = "T_23"
  $ws.Range("B31").Value = "N/A"
  $ws.Range("G31").Value = "N/A"

  $ws.Range("A32").Value = "N_e"
  $ws.Range("B32").Formula = "=B27/4*1000000000"
  $ws.Range("C32:K32").Formula = "=C27/4*1000000000"

  $ws.Range("A33").Value = "N_a"
  $ws.Range("B33").Formula = "=B32/B28"
  $ws.Range("C33:K33").Formula = "=C32/C28"

  $ws.Range("A34").Value = "T_12 (scaled)"
  $ws.Range("B34").Formula = "=2*B30*B33"
  $ws.Range("C34:K34").Formula = "=2*C30*C33"

  $ws.Range("A35").Value = "T_23 (scaled)"
  $ws.Range("B35").Formula = "=2*B31*B33"
  $ws.Range("C35:K35").Formula = "=2*C31*C33"

  $ws.Range("A37").Value = "three_epoch"

  $ws.Range("B38").Value = "b_vulgatus"
  $ws.Range("C38").Value = "b_ovatus"
  $ws.Range("D38").Value = "a_putredinis"
  $ws.Range("E38").Value = "b_uniformis"
  $ws.Range("F38").Value = "e_rectale"
  $ws.Range("G38").Value = "p_copri"
  $ws.Range("H38").Value = "o_splanchnicus"
  $ws.Range("I38").Value = "a_shahii"
  $ws.Range("J38").Value = "coprococcus_sp"

  $ws.Range("A39").Value = "Theta"
  $ws.Range("B39").Value = 2028.0812273500001
  $ws.Range("G39").Value = 6405.0042083099997

  $ws.Range("A40").Value = "Nu_a"
  $ws.Range("B40").Value = 2.64925027
  $ws.Range("G40").Value = 1.69793265

  $ws.Range("A41").Value = "Nu_b"
  $ws.Range("B41").Value = 1.49440752
  $ws.Range("G41").Value = 4.91266537

  $ws.Range("A42").Value = "T_12"
  $ws.Range("B42").Value = 0.00631432
  $ws.Range("G42").Value = 0.46533746999999998

  $ws.Range("A43").Value = "T_23"
  $ws.Range("B43").Value = 0.43254745
  $ws.Range("G43").Value = 0.27208793999999997

  $ws.Range("A44").Value = "N_e"
  $ws.Range("B44").Formula = "=B39/4*1000000000"
  $ws.Range("C44:K44").Formula = "=C39/4*1000000000"

  $ws.Range("A45").Value = "N_a"
  $ws.Range("B45").Formula = "=B44/B40"
  $ws.Range("C45:K45").Formula = "=C44/C40"

  $ws.Range("A46").Value = "T_12 (scaled)"
  $ws.Range("B46").Formula = "=2*B42*B45"
  $ws.Range("C46:K46").Formula = "=2*C42*C45"

  $ws.Range("A47").Value = "T_23 (scaled)"
  $ws.Range("B47").Formula = "=2*B43*B45"
  $ws.Range("C47:K47").Formula = "=2*C43*C45"

# --- scroll/selection state on the new sheet ---------------------------
$ws.Range("A22").Select()
$ws.Range("D27").Select()
